$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the subnet mask text for the /30 WAN rows (D12:D14 share this string):
# "/30`n255.255.255.254" -> "/30`n255.255.255.252"
$ws.Range("D12:D14").Value = "/30" + [char]10 + "255.255.255.252"

# Remove the stray "WAN-1 / next-hop" note block that used to sit next to the
# subnetting table (J6:K7). K7 keeps the same "Output" style that J6/K6 use,
# just emptied out; J7 disappears completely (no value, no style override).
$ws.Range("J6").Copy()
$ws.Range("K7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("J6").ClearContents()
$ws.Range("K6").ClearContents()
$ws.Range("J7").Clear()
$ws.Range("K7").ClearContents()

# Move the active selection like the author did
$ws.Range("F15").Select()
